$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A so the existing "Modelo"/"Fase"/... headers
# and data shift one column to the right, making room for a new "Substância" column.
$ws.Columns.Item(1).Insert()

# Copy the header formatting (bold font, borders, centered alignment) from the
# neighboring header cell onto the newly inserted header cell.
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Set the header row text (A1 is new; B1:E1 already hold the shifted-over values,
# but we (re)write them explicitly to match the target exactly).
$ws.Range("A1").Value = "Substância"
$ws.Range("B1").Value = "Modelo"
$ws.Range("C1").Value = "Fase"
$ws.Range("D1").Value = "Volume Molar (cm³/mol)"
$ws.Range("E1").Value = "Massa Específica (g/cm³)"

# Replace row 2 data with the corrected substance/model/phase/values.
$ws.Range("A2").Value = "Cloreto de Metila"
$ws.Range("B2").Value = "PR"
$ws.Range("C2").Value = "Vapor"
$ws.Range("D2").Value = 2019.1856
$ws.Range("E2").Value = 0.00099

# The old row 3 (second VDW result row) is no longer needed; remove it entirely.
$ws.Rows.Item(3).Delete()
